$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.131.52'
$ws.Range("E2").Value = '  +3.19%  '
$ws.Range("D3").Value = '1.655.69'
$ws.Range("E3").Value = '  +3.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E8").Value = '  +2.09%  '
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("E10").Value = '  +2.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '1.890.99'
$ws.Range("E12").Value = '  +3.53%  '
$ws.Range("D13").Value = '1.649.88'
$ws.Range("E13").Value = '  +2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '241.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.82%  '
$ws.Range("D18").Value = '27.095.40'
$ws.Range("E18").Value = '  +3.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.34%  '
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.28%  '
$ws.Range("E24").Value = '  +3.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("E27").Value = '  +2.32%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0496'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("D32").Value = '1.524.03'
$ws.Range("E32").Value = '  +5.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("E34").Value = '  +2.85%  '
$ws.Range("E35").Value = '  +6.05%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.892'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.43%  '
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.63%  '
$ws.Range("D44").Value = '1.796.85'
$ws.Range("E44").Value = '  +3.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.770'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.911'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.64%  '
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0980'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("E51").Value = '  +1.73%  '
